$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price and volume(1h) data
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "59.668.96"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.60%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.649.07"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.62%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "537.08"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -1.73%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.16"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +3.51%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +1.26%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "6.82"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +5.33%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -0.46%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +1.39%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.12%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.121.69"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +1.79%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "59.581.92"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.55%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "21.40"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +4.14%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.674.42"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +1.93%  "
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.92%  "
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +2.60%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "339.71"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -1.08%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.35"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +2.23%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.20"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -3.25%  "
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.00%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "66.60"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -1.26%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.417"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +2.35%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.55%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.54%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +1.28%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0₃0748"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +1.41%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.07%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.66"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -3.46%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.86"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.96%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "18.89"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.74%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "151.01"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +1.28%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.01"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.72%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +2.59%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.841"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +3.52%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.49%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -0.94%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.61"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "286.25"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.08%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.606"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +1.58%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "10.74"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.00%  "
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +2.83%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "19.30"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +3.72%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.94%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0227"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +1.84%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.966.31"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +1.08%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.56"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +1.14%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "18.40"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.34%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "111.93"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.92%  "
